$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) to reflect the new "through" date
$ws.Name = "Through 2022-12-14"

# Update the column header label (shared string) for the "through" date
$ws.Range("I1").Value = "2022 (through 12-14)"

# Update December carjacking total (2022 column) and yearly Total row
$ws.Range("I13").Value = 60
$ws.Range("I14").Value = 1576
